# LOBSTAHS basic component matrix workbook — add additional notes to the
# "Notes" worksheet changelog/instructions section, and leave the
# workbook with the "Notes" tab active (matching the author's last
# on-screen state before saving).

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")

# Make room for two more instruction rows: rows 9-12 held 4 bullet points;
# the commit expands this to 6 bullet points (rows 9-14), with the old
# row-12 note ("Note: An adduct hierarchy must be specified...") sliding
# down to row 14 unchanged.
$notes.Rows.Item(12).Insert()
$notes.Rows.Item(12).Insert()

# Rewrite the step-by-step instructions (rows 9-13). Order of assignment
# mirrors the order the author appears to have entered/edited them in.
$notes.Range("B10").Value = 'When done with edits/additions, copy all data in "Elemental composition matrix" except for the first row and last column'
$notes.Range("B11").Value = 'Navigate to second tab, "For export to .csv," then paste the data copied from the first worksheet "as values"'
$notes.Range("B9").Value = 'Make any additions or modifications to the default set of entries by editing data in the first worksheet ("Elemental composition matrix")'
$notes.Range("B12").Value = 'Export the "For export to .csv" worksheet as a .csv file with the filename "LOBSTAHS_basic_component_matrix.csv"; this file can be sourced by specifying the file path in the LOBSTAHS database generation function'
$notes.Range("B13").Value = 'Note any changes in changelog (this "Notes" worksheet); consider requesting incorporation of your new lipids into the default LOBSTAHS database via a pull request on GitHub'

# Update on-screen selections for the sheets the author touched, and leave
# the "Notes" tab as the active/selected one when the file was saved.
$compMatrix = $wb.Worksheets.Item(1)
$compMatrix.Activate()
$compMatrix.Range("L23").Select()

$exportCsv = $wb.Worksheets.Item("For export to .csv")
$exportCsv.Activate()
$exportCsv.Range("G55").Select()

$notes.Activate()
$notes.Range("B14").Select()
